# Trade #6 closed at 2026-02-16 21:20:55 - momentum DOWN +0.000%
# Append a new trade row (row 3) to the "momentum" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

$ws.Cells.Item(3, 1).Value = 6
# Force the date-shaped string to stay text instead of being auto-parsed
# into a date serial: a leading quote-prefix keeps Excel's input parser
# from coercing it to a date serial, then ClearFormats() drops the
# resulting quotePrefix style so the cell ends up unstyled (matches the
# other rows).
$ws.Cells.Item(3, 2).Value = "'2026-02-16"
$ws.Cells.Item(3, 2).ClearFormats()
$ws.Cells.Item(3, 3).Value = "21:20:55"
$ws.Cells.Item(3, 4).Value = "momentum"
$ws.Cells.Item(3, 5).Value = "DOWN"
$ws.Cells.Item(3, 6).Value = 69419.005
# Plain Value = "" clears the cell instead of leaving an empty-text cell,
# so use a bare quote-prefix (empty quoted text) then drop the resulting
# quotePrefix style so the cell matches the unstyled empty-string cells
# already on row 2 (G2 / M2).
$ws.Cells.Item(3, 7).Value = "'"
$ws.Cells.Item(3, 7).ClearFormats()
$ws.Cells.Item(3, 8).Value = "OPEN"
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.9
$ws.Cells.Item(3, 12).Value = "Downward momentum: -0.243% over 10 samples"
$ws.Cells.Item(3, 13).Value = "'"
$ws.Cells.Item(3, 13).ClearFormats()
$ws.Cells.Item(3, 14).Value = 0
